# Update WBScanCellAnalysis scan results (re-run scan analysis / plot)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New scan_x, scan_y, scan_z, delta_x, delta_y, delta_z, euclidean_distance values
# for data rows 2-19 (columns F:L), recomputed from the latest scan pass.
$scanData = @(
  @(71.58518981933594, 0.00000002980232238769531, 0, 0.0851898193359375, 0.00000002980232238769531, 0, 0.08518981933594272),
  @(-71.58518981933594, -0.00000002980232238769531, 0, -0.0851898193359375, -0.00000002980232238769531, 0, 0.08518981933594272),
  @(319.8941040039062, 397.6437072753906, 116.9959335327148, -1.60589599609375, 11.79434204101562, -1.892738342285156, 12.05271191736926),
  @(-327.3821411132812, 387.7880554199219, 113.2660598754883, -5.88214111328125, 1.938690185546875, -5.622611999511719, 8.364918972060789),
  @(330.939453125, -386.9641723632812, 112.4130554199219, 9.439453125, -1.11480712890625, -6.475616455078125, 11.50129030617917),
  @(-314.839599609375, -391.5144653320312, 114.5457534790039, 6.660400390625, -5.66510009765625, -4.342918395996094, 9.76295204711074),
  @(175.0008392333984, 2.512862205505371, 120.7578048706055, 1.732269287109375, 2.512862205505371, -4.321296691894531, 5.290447849312793),
  @(-180.8824920654297, -1.754476547241211, 120.9316253662109, -7.613922119140625, -1.754476547241211, -4.147476196289062, 8.845991001009592),
  @(-437.4225924276399, 381.9615940344092, -6.630470556753934, -14.1540233337879, 0.1249851476904382, -0.5044939942539344, 14.16356282899165),
  @(-506.7551674841308, 381.9615940344092, -6.630470556753934, -11.98659839027886, 0.1249851476904382, -0.5044939942539344, 11.99786132805415),
  @(-506.7551674841308, -3.757515036852001, 122.4495703003242, -11.98661322284937, -3.75750740745747, -2.629531262175803, 12.83402482193505),
  @(-428.5609914843379, -392.6309405816843, -2.331833789497878, -5.29242239048591, -10.79433169496554, 3.794142773002122, 12.60646068022973),
  @(-506.7551674841308, -392.6309405816843, -2.331833789497878, -11.98659839027886, -10.79433169496554, 3.794142773002122, 16.57080737600608),
  @(501.180500954079, 394.735955483039, -12.07237651241684, 6.411931860227071, 12.8993465963203, -5.946399949916838, 15.58414852202603),
  @(415.5005639093517, 394.735955483039, -12.07237651241684, -7.768005184500282, 12.8993465963203, -5.946399949916838, 16.18933968767085),
  @(501.180500954079, 11.92230898046362, 115.3838192856772, 6.411946692797585, 11.92231660985816, -9.695282276822837, 16.65092166110491),
  @(501.180500954079, -377.2416097670476, -10.06250876749758, 6.411931860227071, 4.594999119671172, -3.936532204997576, 8.816074687237904),
  @(435.7603124203657, -377.2416097670476, -10.06250876749758, 12.4917433265137, 4.594999119671172, -3.936532204997576, 13.87998393537476)
)

$startRow = 2
$startCol = 6  # column F

for ($r = 0; $r -lt $scanData.Count; $r++) {
    $rowValues = $scanData[$r]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $startCol + $c).Value = $rowValues[$c]
    }
}
